$wb = $excel.ActiveWorkbook

# Add new worksheet "ODI Batting Extra" after the last existing sheet
# (Worksheets.Add(Before, After) - passing After as the current last sheet
# appends the new tab at the end, matching the source workbook's sheet order)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "ODI Batting Extra"

# Match the page margins used on the other sheets (0.75"/0.75"/1"/1" with
# 0.5" header/footer) instead of Excel's narrower tab defaults.
$pageSetup = $newSheet.PageSetup
$pageSetup.LeftMargin = 54
$pageSetup.RightMargin = 54
$pageSetup.TopMargin = 72
$pageSetup.BottomMargin = 72
$pageSetup.HeaderMargin = 36
$pageSetup.FooterMargin = 36

# Header row (column titles)
$newSheet.Range("A1").Value = "MATCH_CODE"
$newSheet.Range("B1").Value = "BATTING_POSITION"
$newSheet.Range("C1").Value = "NUM_4"
$newSheet.Range("D1").Value = "NUM_6"
$newSheet.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$newSheet.Range("F1").Value = "MAN_OF_MATCH"

# Re-use the exact same header style used on the other sheets (bold, bordered,
# centered) by copying an existing header cell's formatting rather than
# re-building it by hand (which would register new style entries).
$headerSource = $wb.Worksheets.Item("ODI Batting").Range("A1")
$headerSource.Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122)

# Data row. MATCH_CODE is copied (as a value-only paste) from the matching
# "ODI Batting" row so the digits are carried over as text (e.g. "4727"),
# the same way the source data is stored, instead of becoming a number.
$matchCodeSource = $wb.Worksheets.Item("ODI Batting").Range("D2")
$matchCodeSource.Copy()
$newSheet.Range("A2").PasteSpecial(-4163)

$newSheet.Range("B2").Value = ""
$newSheet.Range("C2").Value = ""
$newSheet.Range("D2").Value = ""
$newSheet.Range("E2").Value = ""
$newSheet.Range("F2").Value = "NO"

# Restore the original active sheet/selection (adding a sheet makes it the
# active tab, but the workbook's view state otherwise didn't change).
$wb.Worksheets.Item("Player Info").Activate()
